# Rename the worksheet "TPL_rawFilter_metrics" to "TPL_raw".
# (Excel automatically updates dependent refs such as the hidden
#  ExternalData_1 defined name that points at the sheet's range.)
$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "TPL_rawFilter_metrics") {
        $ws = $sheet
        break
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Name = "TPL_raw"
